$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.563.41'
$ws.Range("E2").Value = '  +0.40%  '

$ws.Range("D3").Value = '3.331.01'
$ws.Range("E3").Value = '  +5.84%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '622.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.41'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +25.84%  '

$ws.Range("E8").Value = '  +1.62%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  +9.74%  '

$ws.Range("D11").Value = '3.329.00'
$ws.Range("E11").Value = '  +5.76%  '

$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("B13").Value = 'WrappedBTC'
$ws.Range("C13").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D13").Value = '98.239.68'
$ws.Range("E13").Value = '  +0.68%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.63%  '

$ws.Range("E15").Value = '  +3.14%  '

$ws.Range("D16").Value = '3.937.87'
$ws.Range("E16").Value = '  +5.77%  '

$ws.Range("E17").Value = '  +1.20%  '

$ws.Range("D18").Value = '3.326.83'
$ws.Range("E18").Value = '  +5.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '485.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000209'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.83%  '

$ws.Range("D28").Value = '3.497.92'
$ws.Range("E28").Value = '  +5.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.285'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +19.84%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.195'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +11.98%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.136'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.996'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.76'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.149'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.22'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("E38").Value = '  +3.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '493.70'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.456'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.71'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.86%  '

$ws.Range("E43").Value = '  +3.57%  '

$ws.Range("E44").Value = '  +6.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.778'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +12.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '159.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.00%  '

$ws.Range("E48").Value = '  +1.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.842'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.91%  '
